$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-cell tweaks ---
$ws.Cells.Item(164, 44).Value = 15.0485437   # AR164 (Puerto Rico)
$ws.Cells.Item(199, 7).Value = 17.9978084    # G199 (California)
$ws.Cells.Item(199, 21).Value = 25.4871452   # U199 (Kentucky)

# --- Fill in remaining data for existing row 200 (17 08 2020) ---
$ws.Cells.Item(200, 2).Value = 23.3333333
$ws.Cells.Item(200, 3).Value = 33.3543238
$ws.Cells.Item(200, 4).Value = 27.2249111
$ws.Cells.Item(200, 6).Value = 21.482007
$ws.Cells.Item(200, 7).Value = 18.003383
$ws.Cells.Item(200, 8).Value = 16.7421134
$ws.Cells.Item(200, 9).Value = 9.3447581
$ws.Cells.Item(200, 10).Value = 12.6838235
$ws.Cells.Item(200, 11).Value = 12.1815807
$ws.Cells.Item(200, 12).Value = 24.0047943
$ws.Cells.Item(200, 13).Value = 29.7718869
$ws.Cells.Item(200, 15).Value = 16.9384058
$ws.Cells.Item(200, 16).Value = 26.2112133
$ws.Cells.Item(200, 17).Value = 29.2305812
$ws.Cells.Item(200, 18).Value = 19.1434957
$ws.Cells.Item(200, 19).Value = 24.1831937
$ws.Cells.Item(200, 20).Value = 25.0630372
$ws.Cells.Item(200, 21).Value = 25.634058
$ws.Cells.Item(200, 22).Value = 31.6419783
$ws.Cells.Item(200, 23).Value = 10.7846997
$ws.Cells.Item(200, 24).Value = 13.6341334
$ws.Cells.Item(200, 25).Value = 9.738955799999999
$ws.Cells.Item(200, 26).Value = 16.2143534
$ws.Cells.Item(200, 27).Value = 19.0648143
$ws.Cells.Item(200, 28).Value = 25.8649965
$ws.Cells.Item(200, 30).Value = 33.8170998
$ws.Cells.Item(200, 31).Value = 23.4340164
$ws.Cells.Item(200, 32).Value = 19.5538157
$ws.Cells.Item(200, 33).Value = 24.6112745
$ws.Cells.Item(200, 34).Value = 25.1362022
$ws.Cells.Item(200, 35).Value = 10.432519
$ws.Cells.Item(200, 36).Value = 10.6698762
$ws.Cells.Item(200, 37).Value = 17.6577112
$ws.Cells.Item(200, 38).Value = 23.6675466
$ws.Cells.Item(200, 39).Value = 10.463859
$ws.Cells.Item(200, 40).Value = 19.8511686
$ws.Cells.Item(200, 41).Value = 26.3138364
$ws.Cells.Item(200, 42).Value = 15.1038825
$ws.Cells.Item(200, 43).Value = 14.0953533
$ws.Cells.Item(200, 45).Value = 12.5399813
$ws.Cells.Item(200, 46).Value = 24.9013276
$ws.Cells.Item(200, 47).Value = 23.5124181
$ws.Cells.Item(200, 48).Value = 29.5897621
$ws.Cells.Item(200, 49).Value = 25.4639555
$ws.Cells.Item(200, 50).Value = 23.8650278
$ws.Cells.Item(200, 51).Value = 17.8691015
$ws.Cells.Item(200, 53).Value = 9.5836921
$ws.Cells.Item(200, 54).Value = 17.4154578
$ws.Cells.Item(200, 55).Value = 20.8288989
$ws.Cells.Item(200, 56).Value = 19.7258474
$ws.Cells.Item(200, 57).Value = 18.6876285

# --- New row 201 (18 08 2020) ---
$ws.Cells.Item(201, 1).Value = "18 08 2020"
$ws.Cells.Item(201, 2).Value = 21.4068826
$ws.Cells.Item(201, 3).Value = 31.9193523
$ws.Cells.Item(201, 4).Value = 27.4232695
$ws.Cells.Item(201, 6).Value = 21.1179561
$ws.Cells.Item(201, 7).Value = 17.5810259
$ws.Cells.Item(201, 8).Value = 16.3726729
$ws.Cells.Item(201, 9).Value = 9.437185899999999
$ws.Cells.Item(201, 10).Value = 12.5498008
$ws.Cells.Item(201, 11).Value = 11.9469246
$ws.Cells.Item(201, 12).Value = 23.5529611
$ws.Cells.Item(201, 13).Value = 29.5096159
$ws.Cells.Item(201, 15).Value = 16.9944683
$ws.Cells.Item(201, 16).Value = 26.2761885
$ws.Cells.Item(201, 17).Value = 28.8737433
$ws.Cells.Item(201, 18).Value = 19.3430363
$ws.Cells.Item(201, 19).Value = 24.0310624
$ws.Cells.Item(201, 20).Value = 24.6884001
$ws.Cells.Item(201, 21).Value = 25.9725968
$ws.Cells.Item(201, 22).Value = 31.1415094
$ws.Cells.Item(201, 23).Value = 10.9233726
$ws.Cells.Item(201, 24).Value = 13.8274484
$ws.Cells.Item(201, 25).Value = 9.299959299999999
$ws.Cells.Item(201, 26).Value = 15.8600612
$ws.Cells.Item(201, 27).Value = 19.0632918
$ws.Cells.Item(201, 28).Value = 26.5757572
$ws.Cells.Item(201, 30).Value = 33.3385051
$ws.Cells.Item(201, 31).Value = 22.0340813
$ws.Cells.Item(201, 32).Value = 19.5003739
$ws.Cells.Item(201, 33).Value = 25.7181519
$ws.Cells.Item(201, 34).Value = 25.4531641
$ws.Cells.Item(201, 35).Value = 11.2510675
$ws.Cells.Item(201, 36).Value = 10.5807966
$ws.Cells.Item(201, 37).Value = 16.7992982
$ws.Cells.Item(201, 38).Value = 22.5781586
$ws.Cells.Item(201, 39).Value = 10.4803359
$ws.Cells.Item(201, 40).Value = 19.8338968
$ws.Cells.Item(201, 41).Value = 25.8544159
$ws.Cells.Item(201, 42).Value = 15.0019011
$ws.Cells.Item(201, 43).Value = 13.9073816
$ws.Cells.Item(201, 45).Value = 13.6946141
$ws.Cells.Item(201, 46).Value = 24.1643917
$ws.Cells.Item(201, 47).Value = 24.4186851
$ws.Cells.Item(201, 48).Value = 28.4340552
$ws.Cells.Item(201, 49).Value = 25.1296136
$ws.Cells.Item(201, 50).Value = 23.8942643
$ws.Cells.Item(201, 51).Value = 17.462973
$ws.Cells.Item(201, 53).Value = 8.2478677
$ws.Cells.Item(201, 54).Value = 17.1095745
$ws.Cells.Item(201, 55).Value = 20.7160495
$ws.Cells.Item(201, 56).Value = 20.2573711
$ws.Cells.Item(201, 57).Value = 19.3818245

# --- New row 202 (19 08 2020) ---
$ws.Cells.Item(202, 1).Value = "19 08 2020"
$ws.Cells.Item(202, 2).Value = 23.0020492
$ws.Cells.Item(202, 3).Value = 32.505996
$ws.Cells.Item(202, 4).Value = 27.8943287
$ws.Cells.Item(202, 6).Value = 20.9632856
$ws.Cells.Item(202, 7).Value = 17.4768549
$ws.Cells.Item(202, 8).Value = 16.2517175
$ws.Cells.Item(202, 9).Value = 9.0158796
$ws.Cells.Item(202, 10).Value = 13.2111252
$ws.Cells.Item(202, 11).Value = 12.3198353
$ws.Cells.Item(202, 12).Value = 23.0200266
$ws.Cells.Item(202, 13).Value = 29.1382074
$ws.Cells.Item(202, 15).Value = 18.5314685
$ws.Cells.Item(202, 16).Value = 27.3415708
$ws.Cells.Item(202, 17).Value = 28.4483909
$ws.Cells.Item(202, 18).Value = 19.5452273
$ws.Cells.Item(202, 19).Value = 24.1626899
$ws.Cells.Item(202, 20).Value = 24.421795
$ws.Cells.Item(202, 21).Value = 26.0295216
$ws.Cells.Item(202, 22).Value = 31.3025182
$ws.Cells.Item(202, 23).Value = 10.911696
$ws.Cells.Item(202, 24).Value = 13.711017
$ws.Cells.Item(202, 25).Value = 8.948351600000001
$ws.Cells.Item(202, 26).Value = 15.5806453
$ws.Cells.Item(202, 27).Value = 18.8110252
$ws.Cells.Item(202, 28).Value = 25.9969478
$ws.Cells.Item(202, 30).Value = 33.7816684
$ws.Cells.Item(202, 31).Value = 21.030348
$ws.Cells.Item(202, 32).Value = 19.4797132
$ws.Cells.Item(202, 33).Value = 26.1980023
$ws.Cells.Item(202, 34).Value = 25.3180958
$ws.Cells.Item(202, 35).Value = 10.4273424
$ws.Cells.Item(202, 36).Value = 10.8592962
$ws.Cells.Item(202, 37).Value = 16.8481438
$ws.Cells.Item(202, 38).Value = 22.9262175
$ws.Cells.Item(202, 39).Value = 10.7884875
$ws.Cells.Item(202, 40).Value = 19.7284507
$ws.Cells.Item(202, 41).Value = 26.6187818
$ws.Cells.Item(202, 42).Value = 15.0100735
$ws.Cells.Item(202, 43).Value = 13.7126537
$ws.Cells.Item(202, 45).Value = 13.0053351
$ws.Cells.Item(202, 46).Value = 23.7487206
$ws.Cells.Item(202, 47).Value = 24.8820958
$ws.Cells.Item(202, 48).Value = 28.7104449
$ws.Cells.Item(202, 49).Value = 25.2940884
$ws.Cells.Item(202, 50).Value = 22.8733311
$ws.Cells.Item(202, 51).Value = 17.9279585
$ws.Cells.Item(202, 53).Value = 8.5465711
$ws.Cells.Item(202, 54).Value = 16.9326585
$ws.Cells.Item(202, 55).Value = 20.5022462
$ws.Cells.Item(202, 56).Value = 20.3247245
$ws.Cells.Item(202, 57).Value = 19.1498765

# --- New row 203 (20 08 2020) ---
$ws.Cells.Item(203, 1).Value = "20 08 2020"
$ws.Cells.Item(203, 2).Value = 22.7646454
$ws.Cells.Item(203, 3).Value = 31.7538421
$ws.Cells.Item(203, 4).Value = 26.4174931
$ws.Cells.Item(203, 6).Value = 20.8035111
$ws.Cells.Item(203, 7).Value = 17.2614686
$ws.Cells.Item(203, 8).Value = 16.1052188
$ws.Cells.Item(203, 9).Value = 9.164859
$ws.Cells.Item(203, 10).Value = 13.3821571
$ws.Cells.Item(203, 11).Value = 11.8448637
$ws.Cells.Item(203, 12).Value = 22.8988285
$ws.Cells.Item(203, 13).Value = 28.8857727
$ws.Cells.Item(203, 15).Value = 17.4463938
$ws.Cells.Item(203, 16).Value = 27.8489894
$ws.Cells.Item(203, 17).Value = 28.8394072
$ws.Cells.Item(203, 18).Value = 19.4249427
$ws.Cells.Item(203, 19).Value = 24.4896847
$ws.Cells.Item(203, 20).Value = 24.7260763
$ws.Cells.Item(203, 21).Value = 25.2708323
$ws.Cells.Item(203, 22).Value = 31.3063652
$ws.Cells.Item(203, 23).Value = 10.6267023
$ws.Cells.Item(203, 24).Value = 13.5893541
$ws.Cells.Item(203, 25).Value = 9.217114199999999
$ws.Cells.Item(203, 26).Value = 15.6733471
$ws.Cells.Item(203, 27).Value = 18.7986803
$ws.Cells.Item(203, 28).Value = 26.284885
$ws.Cells.Item(203, 30).Value = 34.7820835
$ws.Cells.Item(203, 31).Value = 22.8202734
$ws.Cells.Item(203, 32).Value = 19.4402246
$ws.Cells.Item(203, 33).Value = 26.9873024
$ws.Cells.Item(203, 34).Value = 25.7048679
$ws.Cells.Item(203, 35).Value = 10.7905498
$ws.Cells.Item(203, 36).Value = 10.6656505
$ws.Cells.Item(203, 37).Value = 16.0517053
$ws.Cells.Item(203, 38).Value = 22.232371
$ws.Cells.Item(203, 39).Value = 10.784958
$ws.Cells.Item(203, 40).Value = 18.8847806
$ws.Cells.Item(203, 41).Value = 26.3027979
$ws.Cells.Item(203, 42).Value = 14.895803
$ws.Cells.Item(203, 43).Value = 13.7752975
$ws.Cells.Item(203, 45).Value = 13.1424154
$ws.Cells.Item(203, 46).Value = 23.8013851
$ws.Cells.Item(203, 47).Value = 25.3049759
$ws.Cells.Item(203, 48).Value = 27.7121714
$ws.Cells.Item(203, 49).Value = 24.9840363
$ws.Cells.Item(203, 50).Value = 22.3643953
$ws.Cells.Item(203, 51).Value = 17.9162484
$ws.Cells.Item(203, 53).Value = 9.113303999999999
$ws.Cells.Item(203, 54).Value = 16.9236381
$ws.Cells.Item(203, 55).Value = 20.5952041
$ws.Cells.Item(203, 56).Value = 19.2633776
$ws.Cells.Item(203, 57).Value = 18.6012902

# --- New row 204 (21 08 2020) ---
$ws.Cells.Item(204, 1).Value = "21 08 2020"
$ws.Cells.Item(204, 2).Value = 23.5070575
$ws.Cells.Item(204, 3).Value = 31.5645889
$ws.Cells.Item(204, 4).Value = 26.6274285
$ws.Cells.Item(204, 6).Value = 20.2979847
$ws.Cells.Item(204, 7).Value = 16.9142816
$ws.Cells.Item(204, 8).Value = 16.5434194
$ws.Cells.Item(204, 9).Value = 9.1456736
$ws.Cells.Item(204, 10).Value = 13.8283379
$ws.Cells.Item(204, 11).Value = 11.5625
$ws.Cells.Item(204, 12).Value = 22.5974427
$ws.Cells.Item(204, 13).Value = 29.1834298
$ws.Cells.Item(204, 15).Value = 17.8869621
$ws.Cells.Item(204, 16).Value = 28.078152
$ws.Cells.Item(204, 17).Value = 29.0149486
$ws.Cells.Item(204, 18).Value = 19.4138426
$ws.Cells.Item(204, 19).Value = 25.4449598
$ws.Cells.Item(204, 20).Value = 25.0000021
$ws.Cells.Item(204, 21).Value = 25.5825675
$ws.Cells.Item(204, 22).Value = 30.424625
$ws.Cells.Item(204, 23).Value = 10.5539775
$ws.Cells.Item(204, 24).Value = 13.7274383
$ws.Cells.Item(204, 25).Value = 9.257622599999999
$ws.Cells.Item(204, 26).Value = 15.2994647
$ws.Cells.Item(204, 27).Value = 18.8595283
$ws.Cells.Item(204, 28).Value = 25.8862777
$ws.Cells.Item(204, 30).Value = 33.5325969
$ws.Cells.Item(204, 31).Value = 22.0088504
$ws.Cells.Item(204, 32).Value = 19.1561488
$ws.Cells.Item(204, 33).Value = 25.630872
$ws.Cells.Item(204, 34).Value = 25.9058954
$ws.Cells.Item(204, 35).Value = 11.103664
$ws.Cells.Item(204, 36).Value = 10.7439521
$ws.Cells.Item(204, 37).Value = 16.9772481
$ws.Cells.Item(204, 38).Value = 22.4373924
$ws.Cells.Item(204, 39).Value = 10.6688272
$ws.Cells.Item(204, 40).Value = 19.2142415
$ws.Cells.Item(204, 41).Value = 25.5380217
$ws.Cells.Item(204, 42).Value = 15.1521182
$ws.Cells.Item(204, 43).Value = 13.8658093
$ws.Cells.Item(204, 45).Value = 12.467359
$ws.Cells.Item(204, 46).Value = 23.5771635
$ws.Cells.Item(204, 47).Value = 23.1843811
$ws.Cells.Item(204, 48).Value = 27.8251686
$ws.Cells.Item(204, 49).Value = 25.0141461
$ws.Cells.Item(204, 50).Value = 22.8090572
$ws.Cells.Item(204, 51).Value = 17.9472948
$ws.Cells.Item(204, 53).Value = 8.805975200000001
$ws.Cells.Item(204, 54).Value = 16.4454095
$ws.Cells.Item(204, 55).Value = 20.5554814
$ws.Cells.Item(204, 56).Value = 20.277054
$ws.Cells.Item(204, 57).Value = 17.8621282

# --- New row 205 (22 08 2020) ---
$ws.Cells.Item(205, 1).Value = "22 08 2020"
$ws.Cells.Item(205, 2).Value = 22.9452055
$ws.Cells.Item(205, 3).Value = 30.8327891
$ws.Cells.Item(205, 4).Value = 25.7999306
$ws.Cells.Item(205, 6).Value = 20.1438353
$ws.Cells.Item(205, 7).Value = 16.5872906
$ws.Cells.Item(205, 8).Value = 16.2995906
$ws.Cells.Item(205, 9).Value = 9.0683632
$ws.Cells.Item(205, 10).Value = 14.4366197
$ws.Cells.Item(205, 11).Value = 11.6932686
$ws.Cells.Item(205, 12).Value = 22.4488329
$ws.Cells.Item(205, 13).Value = 29.0122706
$ws.Cells.Item(205, 15).Value = 17.5241158
$ws.Cells.Item(205, 16).Value = 28.5529128
$ws.Cells.Item(205, 17).Value = 29.2008244
$ws.Cells.Item(205, 18).Value = 19.3972542
$ws.Cells.Item(205, 19).Value = 25.2040445
$ws.Cells.Item(205, 20).Value = 24.7459786
$ws.Cells.Item(205, 21).Value = 25.4570351
$ws.Cells.Item(205, 22).Value = 29.8786055
$ws.Cells.Item(205, 23).Value = 10.3792713
$ws.Cells.Item(205, 24).Value = 13.2965288
$ws.Cells.Item(205, 25).Value = 9.380269800000001
$ws.Cells.Item(205, 26).Value = 15.3897299
$ws.Cells.Item(205, 27).Value = 18.5204141
$ws.Cells.Item(205, 28).Value = 26.6598015
$ws.Cells.Item(205, 30).Value = 33.1915918
$ws.Cells.Item(205, 31).Value = 22.4049036
$ws.Cells.Item(205, 32).Value = 19.216757
$ws.Cells.Item(205, 33).Value = 24.6281559
$ws.Cells.Item(205, 34).Value = 26.6383221
$ws.Cells.Item(205, 35).Value = 11.5602837
$ws.Cells.Item(205, 36).Value = 10.7805547
$ws.Cells.Item(205, 37).Value = 17.4083596
$ws.Cells.Item(205, 38).Value = 22.0098618
$ws.Cells.Item(205, 39).Value = 10.4031549
$ws.Cells.Item(205, 40).Value = 18.8264815
$ws.Cells.Item(205, 41).Value = 26.6788252
$ws.Cells.Item(205, 42).Value = 15.0409724
$ws.Cells.Item(205, 43).Value = 13.7530155
$ws.Cells.Item(205, 45).Value = 12.1960275
$ws.Cells.Item(205, 46).Value = 23.1378347
$ws.Cells.Item(205, 47).Value = 24.7266831
$ws.Cells.Item(205, 48).Value = 27.8562097
$ws.Cells.Item(205, 49).Value = 24.618424
$ws.Cells.Item(205, 50).Value = 22.6705241
$ws.Cells.Item(205, 51).Value = 17.8096348
$ws.Cells.Item(205, 53).Value = 10.461683
$ws.Cells.Item(205, 54).Value = 16.6870979
$ws.Cells.Item(205, 55).Value = 20.3445259
$ws.Cells.Item(205, 56).Value = 20.4930747
$ws.Cells.Item(205, 57).Value = 17.782641

# --- New row 206 (23 08 2020) ---
$ws.Cells.Item(206, 1).Value = "23 08 2020"
$ws.Cells.Item(206, 2).Value = 22.3214286
$ws.Cells.Item(206, 3).Value = 30.8178159
$ws.Cells.Item(206, 4).Value = 25.630216
$ws.Cells.Item(206, 6).Value = 19.7502937
$ws.Cells.Item(206, 7).Value = 16.5269201
$ws.Cells.Item(206, 8).Value = 16.2560963
$ws.Cells.Item(206, 9).Value = 9.648139799999999
$ws.Cells.Item(206, 10).Value = 14.9782923
$ws.Cells.Item(206, 11).Value = 11.3636364
$ws.Cells.Item(206, 12).Value = 22.0998004
$ws.Cells.Item(206, 13).Value = 28.4578547
$ws.Cells.Item(206, 15).Value = 18.5515873
$ws.Cells.Item(206, 16).Value = 29.1333218
$ws.Cells.Item(206, 17).Value = 28.7968408
$ws.Cells.Item(206, 18).Value = 18.9785118
$ws.Cells.Item(206, 19).Value = 25.3275417
$ws.Cells.Item(206, 20).Value = 25.055695
$ws.Cells.Item(206, 21).Value = 25.2063575
$ws.Cells.Item(206, 22).Value = 29.8759385
$ws.Cells.Item(206, 23).Value = 10.6966489
$ws.Cells.Item(206, 24).Value = 13.5102246
$ws.Cells.Item(206, 25).Value = 9.158083299999999
$ws.Cells.Item(206, 26).Value = 15.5759924
$ws.Cells.Item(206, 27).Value = 18.6148504
$ws.Cells.Item(206, 28).Value = 25.6734863
$ws.Cells.Item(206, 30).Value = 33.4500077
$ws.Cells.Item(206, 31).Value = 22.6427579
$ws.Cells.Item(206, 32).Value = 19.4914461
$ws.Cells.Item(206, 33).Value = 25.2137565
$ws.Cells.Item(206, 34).Value = 26.5805954
$ws.Cells.Item(206, 35).Value = 10.5396679
$ws.Cells.Item(206, 36).Value = 11.2534041
$ws.Cells.Item(206, 37).Value = 17.4642652
$ws.Cells.Item(206, 38).Value = 21.6700614
$ws.Cells.Item(206, 39).Value = 10.4719276
$ws.Cells.Item(206, 40).Value = 18.4499495
$ws.Cells.Item(206, 41).Value = 26.780914
$ws.Cells.Item(206, 42).Value = 14.6117776
$ws.Cells.Item(206, 43).Value = 13.9447743
$ws.Cells.Item(206, 45).Value = 11.7337082
$ws.Cells.Item(206, 46).Value = 23.6292604
$ws.Cells.Item(206, 47).Value = 23.8106786
$ws.Cells.Item(206, 48).Value = 27.5000761
$ws.Cells.Item(206, 49).Value = 24.4238878
$ws.Cells.Item(206, 50).Value = 23.4138733
$ws.Cells.Item(206, 51).Value = 17.7422052
$ws.Cells.Item(206, 53).Value = 8.994971400000001
$ws.Cells.Item(206, 54).Value = 16.4278659
$ws.Cells.Item(206, 55).Value = 20.1877464
$ws.Cells.Item(206, 56).Value = 20.0823243
$ws.Cells.Item(206, 57).Value = 17.9415404

# --- New rows 207-208: date labels only, no numeric data ---
$ws.Cells.Item(207, 1).Value = "24 08 2020"
$ws.Cells.Item(208, 1).Value = "25 08 2020"

